$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6
$ws.Range("A6").Value = 9399.64
$ws.Range("B6").Value = 10438.24
$ws.Range("C6").Value = 22.31
$ws.Range("D6").Value = 20.09
$ws.Range("E6").Value = $false
$ws.Range("F6").Value = -9.95
$ws.Range("G6").Value = 42607.884293981479
$ws.Range("H6").Value = $false

# Row 7
$ws.Range("A7").Value = 9918.5
$ws.Range("B7").Value = 9399.64
$ws.Range("C7").Value = 20.3
$ws.Range("D7").Value = 19.18
$ws.Range("E7").Value = $true
$ws.Range("F7").Value = -5.52
$ws.Range("G7").Value = 42608.616365740738
$ws.Range("H7").Value = $true
